# Generate Report for Archive
# - Status moves from "Ready for handoff" to "In Translation" across the
#   Overview summary sheet (columns E/F) and the per-locale detail sheets
#   (zh-cn / de-de, column C "Status").
# - The now-shorter status text lets the report narrow the Status columns,
#   so bring their widths down to match the regenerated report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text ---------------------------------------------------
$overview.Range("E2:E4").Value = "In Translation"
$overview.Range("F2:F4").Value = "In Translation"

$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- Narrow the Status columns to fit the new, shorter text ---------------
$overview.Columns.Item(5).ColumnWidth = 12.43   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.43   # column F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = 12.43       # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.43       # column C (Status)
